# Update results: anonymize "fedcore" -> "approach" and add separator
# borders on the merged-header spacer cells (quality_comparison +
# computational_comparison sheets).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- quality_comparison ------------------------------------------------

# C1: clear style back to default, then give it a top+bottom thin border
# (built as "box all sides" then strip left & right, so the engine reuses
# / creates a single compact style rather than one per edge).
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.LineStyle = 1
$c1.Borders.Item(7).LineStyle = -4142
$c1.Borders.Item(10).LineStyle = -4142

# D1: top+bottom+right thin border (box all sides, strip left only).
$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.LineStyle = 1
$d1.Borders.Item(7).LineStyle = -4142

# Anonymize the header label.
$ws1.Range("C2").Value = "approach"

# --- computational_comparison -------------------------------------------

# Reuse the exact formats just built on sheet 1 (copy/paste-format keeps
# everything on the same two style slots instead of minting new ones).
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$c1.Copy()
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$d1.Copy()
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Anonymize both header labels.
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell at G5.
$ws2.Range("G5").ClearContents()

Write-Host "edits applied"
